$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2959.3333
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H44").Value = 40000
$ws.Range("J44").Value = 40000
$ws.Range("L44").Value = 40000
$ws.Range("N44").Value = -40924
$ws.Range("H96").Value = 3299.4
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 3299.4
$ws.Range("K96").Value = 0
$ws.Range("L96").ClearContents()
$ws.Range("M96").Value = 9898.200000000001
$ws.Range("N96").Value = -12644.2
$ws.Range("H107").Value = 1326.2222
$ws.Range("I107").Value = 988.3333
$ws.Range("K107").Value = 988.3333
$ws.Range("M107").Value = 931.6667
$ws.Range("H116").Value = 11649.368
$ws.Range("I116").Value = 4548.5
$ws.Range("J116").Value = 13542.934
$ws.Range("K116").Value = 4548.5
$ws.Range("L116").Value = 13542.934
$ws.Range("M116").Value = -1106.5
$ws.Range("N116").Value = -20426.934
$ws.Range("H118").Value = 716.3333
$ws.Range("I118").Value = 574.5
$ws.Range("K118").Value = 1723.5
$ws.Range("M118").Value = -66.5
$ws.Range("H137").Value = 3292.077
$ws.Range("J137").Value = 4528.7144
$ws.Range("L137").Value = 13586.1432
$ws.Range("N137").Value = -18686.1432
$ws.Range("H138").Value = 2502.827
$ws.Range("I138").Value = 2048.5144
$ws.Range("J138").Value = 3438.1765
$ws.Range("K138").Value = 6145.5432
$ws.Range("L138").Value = 10314.5295
$ws.Range("M138").Value = -1005.5432
$ws.Range("N138").Value = -20594.5295

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3439.7273
$ws.Range("I2").Value = 3483.8
$ws.Range("J2").Value = 2999
$ws.Range("K2").Value = 3483.8
$ws.Range("L2").Value = 2999
$ws.Range("M2").Value = -3370.8
$ws.Range("N2").Value = -3225
$ws.Range("H74").Value = 3026.3901
$ws.Range("I74").Value = 2961.3784
$ws.Range("K74").Value = 2961.3784
$ws.Range("M74").Value = -2087.3784
$ws.Range("H77").Value = 3026.3901
$ws.Range("I77").Value = 2961.3784
$ws.Range("K77").Value = 14806.892
$ws.Range("M77").Value = -10438.892
$ws.Range("H88").Value = 3516.077
$ws.Range("I88").Value = 2245
$ws.Range("J88").Value = 3747.182
$ws.Range("K88").Value = 2245
$ws.Range("L88").Value = 3747.182
$ws.Range("M88").Value = -1839
$ws.Range("N88").Value = -4559.182
$ws.Range("H91").Value = 3516.077
$ws.Range("I91").Value = 2245
$ws.Range("J91").Value = 3747.182
$ws.Range("K91").Value = 2245
$ws.Range("L91").Value = 3747.182
$ws.Range("M91").Value = -841
$ws.Range("N91").Value = -6555.182
$ws.Range("H102").Value = 3014.5454
$ws.Range("I102").Value = 2782.7
$ws.Range("K102").Value = 2782.7
$ws.Range("M102").Value = -1160.7
$ws.Range("H116").Value = 3439.7273
$ws.Range("I116").Value = 3483.8
$ws.Range("J116").Value = 2999
$ws.Range("K116").Value = 3483.8
$ws.Range("L116").Value = 2999
$ws.Range("M116").Value = -1189.8
$ws.Range("N116").Value = -7587
$ws.Range("H132").Value = 51116.047
$ws.Range("I132").Value = 51116.047
$ws.Range("K132").Value = 153348.141
$ws.Range("M132").Value = -150818.141

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3439.7273
$ws.Range("I3").Value = 3483.8
$ws.Range("J3").Value = 2999
$ws.Range("K3").Value = 3483.8
$ws.Range("L3").Value = 2999
$ws.Range("M3").Value = -3369.8
$ws.Range("N3").Value = -3227
$ws.Range("H94").Value = 3480.6924
$ws.Range("I94").Value = 3225.9
$ws.Range("J94").Value = 4330
$ws.Range("K94").Value = 3225.9
$ws.Range("L94").Value = 4330
$ws.Range("M94").Value = -2774.9
$ws.Range("N94").Value = -5232
$ws.Range("H134").Value = 6597.8
$ws.Range("I134").Value = 5995
$ws.Range("K134").Value = 17985
$ws.Range("M134").Value = -15450

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4505.5
$ws.Range("J31").Value = 5794
$ws.Range("L31").Value = 5794
$ws.Range("N31").Value = -6384
$ws.Range("H34").Value = 4505.5
$ws.Range("J34").Value = 5794
$ws.Range("L34").Value = 5794
$ws.Range("N34").Value = -6198
$ws.Range("H58").Value = 103031.5
$ws.Range("I58").Value = 169135.17
$ws.Range("K58").Value = 169135.17
$ws.Range("M58").Value = -168932.17
$ws.Range("H120").Value = 54999
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 54999
$ws.Range("K120").Value = 0
$ws.Range("L120").ClearContents()
$ws.Range("M120").Value = 54999
$ws.Range("N120").Value = -62257
$ws.Range("H132").Value = 1747.625
$ws.Range("I132").Value = 1712
$ws.Range("J132").Value = 1997
$ws.Range("K132").Value = 5136
$ws.Range("L132").Value = 5991
$ws.Range("M132").Value = -2606
$ws.Range("N132").Value = -11051
$ws.Range("H134").Value = 97773
$ws.Range("I134").Value = 116109.89
$ws.Range("J134").Value = 15257
$ws.Range("K134").Value = 348329.67
$ws.Range("L134").Value = 45771
$ws.Range("M134").Value = -345794.67
$ws.Range("N134").Value = -50841
$ws.Range("H136").Value = 103031.5
$ws.Range("I136").Value = 169135.17
$ws.Range("K136").Value = 507405.51
$ws.Range("M136").Value = -504855.51

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2480.8
$ws.Range("I140").Value = 2095.353
$ws.Range("K140").Value = 6286.059
$ws.Range("M140").Value = -1106.059

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 230163.78
$ws.Range("J113").Value = 340166.66
$ws.Range("L113").Value = 340166.66
$ws.Range("N113").Value = -344506.66
$ws.Range("H122").Value = 3220.2144
$ws.Range("I122").Value = 2251.5293
$ws.Range("K122").Value = 6754.5879
$ws.Range("M122").Value = -4304.5879
$ws.Range("H132").Value = 85848.25
$ws.Range("I132").Value = 113520
$ws.Range("K132").Value = 340560
$ws.Range("M132").Value = -338030

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3195.9473
$ws.Range("I61").Value = 2901.2778
$ws.Range("K61").Value = 2901.2778
$ws.Range("M61").Value = -2699.2778
$ws.Range("H82").Value = 2356.3704
$ws.Range("I82").Value = 2288.077
$ws.Range("J82").Value = 2419.7856
$ws.Range("K82").Value = 2288.077
$ws.Range("L82").Value = 2419.7856
$ws.Range("M82").Value = -1927.077
$ws.Range("N82").Value = -3141.7856
$ws.Range("H85").Value = 2356.3704
$ws.Range("I85").Value = 2288.077
$ws.Range("J85").Value = 2419.7856
$ws.Range("K85").Value = 2288.077
$ws.Range("L85").Value = 2419.7856
$ws.Range("M85").Value = -1040.077
$ws.Range("N85").Value = -4915.7856
$ws.Range("H113").Value = 3195.9473
$ws.Range("I113").Value = 2901.2778
$ws.Range("K113").Value = 2901.2778
$ws.Range("M113").Value = -731.2777999999998
$ws.Range("H133").Value = 88338.836
$ws.Range("J133").Value = 88338.836
$ws.Range("L133").Value = 88338.836
$ws.Range("N133").Value = -93398.836

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 299.66666
$ws.Range("I2").Value = 299
$ws.Range("K2").Value = 299
$ws.Range("M2").Value = -187
$ws.Range("H41").Value = 14365.5
$ws.Range("J41").Value = 14365.5
$ws.Range("L41").Value = 14365.5
$ws.Range("N41").Value = -15145.5
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").ClearContents()
$ws.Range("N102").Value = 0
$ws.Range("H103").Value = 22401.334
$ws.Range("J103").Value = 22401.334
$ws.Range("L103").Value = 22401.334
$ws.Range("N103").Value = -24745.334
$ws.Range("H104").Value = 23000
$ws.Range("J104").Value = 23000
$ws.Range("L104").Value = 23000
$ws.Range("N104").Value = -29988
$ws.Range("H106").Value = 59999.5
$ws.Range("J106").Value = 99999
$ws.Range("L106").Value = 99999
$ws.Range("N106").Value = -102523
$ws.Range("H127").Value = 97428.5
$ws.Range("J127").Value = 97428.5
$ws.Range("L127").Value = 97428.5
$ws.Range("N127").Value = -107348.5
$ws.Range("H133").Value = 90000
$ws.Range("J133").Value = 90000
$ws.Range("L133").Value = 90000
$ws.Range("N133").Value = -100120
